$d = $word.ActiveDocument

# Append " (Changed main)" after the first sentence, split across three
# separate runs (" (", "Changed main", ")") to mirror the target markup.
$r = $d.Content
$r.Find.Execute("This is a Microsoft word document.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(" (Changed main)")

# The insertion above lands in the same run as the preceding sentence
# (identical formatting gets coalesced on save). Force run boundaries at
# each of the three new text chunks by briefly bracketing them with a
# bookmark and removing it again - this splits the run without leaving
# any formatting residue behind.
$body = $d.Content

$r1 = $body.Duplicate
$r1.Find.Execute(" (", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("tmpSplit1", $r1)
$d.Bookmarks("tmpSplit1").Delete()

$r2 = $body.Duplicate
$r2.Find.Execute("Changed main", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("tmpSplit2", $r2)
$d.Bookmarks("tmpSplit2").Delete()

$r3 = $body.Duplicate
$r3.Find.Execute(")", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("tmpSplit3", $r3)
$d.Bookmarks("tmpSplit3").Delete()
